# ADI-scrum-BrndownChart.xlsx edit script
# Commit message: "Categories in DB changed to ENUM"
#
# Reflects splitting the former "Revidering av coden" / "Ändra i API från
# boolean till Enum" task on the "Sprint 2" sheet into five separate tasks
# (one per API: User, Category, Cart, Product, plus two placeholder tasks
# 5.6/5.7), adjusting the burndown sums, and nudging chart/view state that
# moves along with the sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Sprint 2"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sprint 2")

# Insert 5 new rows right after row 15 (old rows 16/17 shift to 21/22).
$ws2.Rows("16:20").Insert()

# Row 15: the original "boolean -> Enum" task becomes specifically about the
# User API, and its remaining estimate drops from 15h to 6h (spread across
# the now-5 split tasks).
$ws2.Range("A15:G15").Rows.RowHeight = 15
$ws2.Range("B15").Value = "Revidering av coden API från boolean till Enum"
$ws2.Range("C15").Value = "Task 5.1"
$ws2.Range("D15").Value = "Ändra User API"
$ws2.Range("E15").Value = 6
$ws2.Range("F15").Value = 6
$ws2.Range("G15").Value = 6

# New row 16: Category API
$ws2.Range("A16").Value = "ADI floggit"
$ws2.Range("B16").Value = "Revidering av coden API från boolean till Enum"
$ws2.Range("C16").Value = "Task 5.2"
$ws2.Range("D16").Value = "Ändra Category API"
$ws2.Range("E16").Value = 2
$ws2.Range("F16").Value = 2
$ws2.Range("G16").Value = 2

# New row 17: Cart Api
$ws2.Range("A17").Value = "ADI floggit"
$ws2.Range("B17").Value = "Revidering av coden API från boolean till Enum"
$ws2.Range("C17").Value = "Task 5.3"
$ws2.Range("D17").Value = "Ändra Cart Api"
$ws2.Range("E17").Value = 2
$ws2.Range("F17").Value = 2
$ws2.Range("G17").Value = 2

# New row 18: Product API
$ws2.Range("A18").Value = "ADI floggit"
$ws2.Range("B18").Value = "Revidering av coden API från boolean till Enum"
$ws2.Range("C18").Value = "Task 5.4"
$ws2.Range("D18").Value = "Ändra Product API "
$ws2.Range("E18").Value = 2
$ws2.Range("F18").Value = 2
$ws2.Range("G18").Value = 2
$ws2.Range("A18:G18").Rows.RowHeight = 15

# New row 19: Task 5.6 (placeholder, no "Vad"/hours yet)
$ws2.Range("A19").Value = "ADI floggit"
$ws2.Range("B19").Value = "Revidering av coden API från boolean till Enum"
$ws2.Range("C19").Value = "Task 5.6"

# New row 20: Task 5.7 (placeholder, no "Vad"/hours yet)
$ws2.Range("A20").Value = "ADI floggit"
$ws2.Range("B20").Value = "Revidering av coden API från boolean till Enum"
$ws2.Range("C20").Value = "Task 5.7"

# Row 21 ("Actual Burndown" summary, shifted down from 16): sums now cover
# the extended task list (rows 3:20 / 14:20 / 15:20).
$ws2.Range("E21").Formula = "=SUM(E3:E20)"
$ws2.Range("F21").Formula = "=SUM(F14:F20)"
$ws2.Range("G21").Formula = "=SUM(G15:G20)"

# Column B needs to be wider to fit the longer task description.
$ws2.Columns("B").ColumnWidth = 36.1640625

# Restore the view's last-selected cell like the authored change.
$ws2.Range("F33").Select()

# ---------------------------------------------------------------------
# Sheet "Sprint 1" -- only the scrolled/selected view position changed.
# ---------------------------------------------------------------------
$ws1b = $wb.Worksheets.Item("Sprint 1")
$ws1b.Application.Goto($ws1b.Range("A8"), $false)

# ---------------------------------------------------------------------
# Sheet "ADI-burndown" -- the burndown row for the second sprint was
# retyped, which re-splits the shared formula group in row 39 and nudges
# the saved scroll/selection position.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ADI-burndown")
$ws1.Range("G39:L39").Formula = "=F39-15"
$ws1.Range("M39:V39").Formula = "=L39-15"

# ---------------------------------------------------------------------
# Chart "chart3" (the Sprint 2 burndown chart) tracks the rows that moved
# from 16/17 to 21/22.
# ---------------------------------------------------------------------
$chartSheet = $wb.Worksheets.Item("Sprint 2")
$chartObj = $chartSheet.ChartObjects(1)
$chart = $chartObj.Chart
$chart.SeriesCollection(1).Formula = "=SERIES(,,'Sprint 2'!`$E`$21:`$O`$21,1)"
$chart.SeriesCollection(2).Formula = "=SERIES(,,'Sprint 2'!`$E`$22:`$O`$22,2)"

# The chart anchor moves down by 5 rows along with the inserted rows.
$chartObj.Top = $chartSheet.Rows("26").Top
$chartObj.Left = $chartSheet.Columns("A").Left

# ---------------------------------------------------------------------
# Restore the final active sheet/selection to match the authored file
# (Sprint 2 stays the tab that was active when the file was saved).
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B79:E84").Select()
